# Apply scheduled-runner updates to the Leve profit tables (columns H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 61.1  # H4
$ws.Cells.Item(4, 9).Value = 51.375  # I4
$ws.Cells.Item(4, 10).Value = 100  # J4
$ws.Cells.Item(4, 11).Value = 51.375  # K4
$ws.Cells.Item(4, 12).Value = 100  # L4
$ws.Cells.Item(4, 13).Value = 62.625  # M4
$ws.Cells.Item(4, 14).Value = -328  # N4

$ws.Cells.Item(12, 8).Value = 282.55554  # H12
$ws.Cells.Item(12, 9).Value = 271.83334  # I12
$ws.Cells.Item(12, 10).Value = 304  # J12
$ws.Cells.Item(12, 11).Value = 271.83334  # K12
$ws.Cells.Item(12, 12).Value = 304  # L12
$ws.Cells.Item(12, 13).Value = -101.83334  # M12
$ws.Cells.Item(12, 14).Value = -644  # N12

$ws.Cells.Item(19, 8).Value = 2342.84  # H19
$ws.Cells.Item(19, 9).Value = 3249.2354  # I19
$ws.Cells.Item(19, 10).Value = 416.75  # J19
$ws.Cells.Item(19, 11).Value = 3249.2354  # K19
$ws.Cells.Item(19, 12).Value = 416.75  # L19
$ws.Cells.Item(19, 13).Value = -3074.2354  # M19
$ws.Cells.Item(19, 14).Value = -766.75  # N19

$ws.Cells.Item(33, 8).Value = 17242216  # H33
$ws.Cells.Item(33, 9).Value = 29412422  # I33
$ws.Cells.Item(33, 11).Value = 29412422  # K33
$ws.Cells.Item(33, 13).Value = -29412193  # M33

$ws.Cells.Item(54, 8).Value = 0  # H54
$ws.Cells.Item(54, 9).Value = 0  # I54
$ws.Cells.Item(54, 11).Value = 0  # K54
$ws.Cells.Item(54, 13).ClearContents()  # M54 removed

$ws.Cells.Item(61, 8).Value = 179.75  # H61
$ws.Cells.Item(61, 9).Value = 179.75  # I61
$ws.Cells.Item(61, 11).Value = 539.25  # K61
$ws.Cells.Item(61, 13).Value = -367.25  # M61

$ws.Cells.Item(86, 8).Value = 3537.558  # H86
$ws.Cells.Item(86, 9).Value = 2142.5881  # I86
$ws.Cells.Item(86, 10).Value = 4449.654  # J86
$ws.Cells.Item(86, 11).Value = 2142.5881  # K86
$ws.Cells.Item(86, 12).Value = 4449.654  # L86
$ws.Cells.Item(86, 13).Value = -1019.5881  # M86
$ws.Cells.Item(86, 14).Value = -6695.654  # N86

$ws.Cells.Item(88, 8).Value = 9414.286  # H88
$ws.Cells.Item(88, 9).Value = 1333.3334  # I88
$ws.Cells.Item(88, 10).Value = 15475  # J88
$ws.Cells.Item(88, 11).Value = 1333.3334  # K88
$ws.Cells.Item(88, 12).Value = 15475  # L88
$ws.Cells.Item(88, 13).Value = -927.3334  # M88
$ws.Cells.Item(88, 14).Value = -16287  # N88

$ws.Cells.Item(89, 8).Value = 3537.558  # H89
$ws.Cells.Item(89, 9).Value = 2142.5881  # I89
$ws.Cells.Item(89, 10).Value = 4449.654  # J89
$ws.Cells.Item(89, 11).Value = 10712.9405  # K89
$ws.Cells.Item(89, 12).Value = 22248.27  # L89
$ws.Cells.Item(89, 13).Value = -5096.940500000001  # M89
$ws.Cells.Item(89, 14).Value = -33480.27  # N89

$ws.Cells.Item(91, 8).Value = 9414.286  # H91
$ws.Cells.Item(91, 9).Value = 1333.3334  # I91
$ws.Cells.Item(91, 10).Value = 15475  # J91
$ws.Cells.Item(91, 11).Value = 1333.3334  # K91
$ws.Cells.Item(91, 12).Value = 15475  # L91
$ws.Cells.Item(91, 13).Value = 70.66660000000002  # M91
$ws.Cells.Item(91, 14).Value = -18283  # N91

$ws.Cells.Item(106, 8).Value = 3551.1538  # H106
$ws.Cells.Item(106, 9).Value = 1694.1666  # I106
$ws.Cells.Item(106, 11).Value = 1694.1666  # K106
$ws.Cells.Item(106, 13).Value = -1063.1666  # M106

$ws.Cells.Item(116, 8).Value = 111363.42  # H116
$ws.Cells.Item(116, 9).Value = 149707.5  # I116
$ws.Cells.Item(116, 11).Value = 149707.5  # K116
$ws.Cells.Item(116, 13).Value = -146265.5  # M116

$ws.Cells.Item(129, 8).Value = 743.9286  # H129
$ws.Cells.Item(129, 9).Value = 354.36365  # I129
$ws.Cells.Item(129, 10).Value = 2172.3333  # J129
$ws.Cells.Item(129, 11).Value = 1063.09095  # K129
$ws.Cells.Item(129, 12).Value = 6516.999899999999  # L129
$ws.Cells.Item(129, 13).Value = 3936.90905  # M129
$ws.Cells.Item(129, 14).Value = -16516.9999  # N129

$ws.Cells.Item(137, 8).Value = 33017.207  # H137
$ws.Cells.Item(137, 9).Value = 42180.31  # I137
$ws.Cells.Item(137, 10).Value = 3237.125  # J137
$ws.Cells.Item(137, 11).Value = 126540.93  # K137
$ws.Cells.Item(137, 12).Value = 9711.375  # L137
$ws.Cells.Item(137, 13).Value = -123990.93  # M137
$ws.Cells.Item(137, 14).Value = -14811.375  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7591.915  # H32
$ws.Cells.Item(32, 9).Value = 6142.0586  # I32
$ws.Cells.Item(32, 10).Value = 16834.75  # J32
$ws.Cells.Item(32, 11).Value = 6142.0586  # K32
$ws.Cells.Item(32, 12).Value = 16834.75  # L32
$ws.Cells.Item(32, 13).Value = -5855.0586  # M32
$ws.Cells.Item(32, 14).Value = -17408.75  # N32

$ws.Cells.Item(45, 8).Value = 1995  # H45
$ws.Cells.Item(45, 9).Value = 1490  # I45
$ws.Cells.Item(45, 10).Value = 2500  # J45
$ws.Cells.Item(45, 11).Value = 1490  # K45
$ws.Cells.Item(45, 12).Value = 2500  # L45
$ws.Cells.Item(45, 13).Value = -1113  # M45
$ws.Cells.Item(45, 14).Value = -3254  # N45

$ws.Cells.Item(61, 8).Value = 3095.365  # H61
$ws.Cells.Item(61, 9).Value = 2246.8276  # I61
$ws.Cells.Item(61, 10).Value = 3819.1177  # J61
$ws.Cells.Item(61, 11).Value = 2246.8276  # K61
$ws.Cells.Item(61, 12).Value = 3819.1177  # L61
$ws.Cells.Item(61, 13).Value = -2034.8276  # M61
$ws.Cells.Item(61, 14).Value = -4243.1177  # N61

$ws.Cells.Item(63, 8).Value = 17024  # H63
$ws.Cells.Item(63, 9).Value = 30000  # I63
$ws.Cells.Item(63, 10).Value = 4048  # J63
$ws.Cells.Item(63, 11).Value = 30000  # K63
$ws.Cells.Item(63, 12).Value = 4048  # L63
$ws.Cells.Item(63, 13).Value = -29314  # M63
$ws.Cells.Item(63, 14).Value = -5420  # N63

$ws.Cells.Item(66, 8).Value = 17024  # H66
$ws.Cells.Item(66, 9).Value = 30000  # I66
$ws.Cells.Item(66, 10).Value = 4048  # J66
$ws.Cells.Item(66, 11).Value = 150000  # K66
$ws.Cells.Item(66, 12).Value = 20240  # L66
$ws.Cells.Item(66, 13).Value = -146568  # M66
$ws.Cells.Item(66, 14).Value = -27104  # N66

$ws.Cells.Item(113, 8).Value = 31603  # H113
$ws.Cells.Item(113, 10).Value = 31603  # J113
$ws.Cells.Item(113, 12).Value = 31603  # L113
$ws.Cells.Item(113, 14).Value = -40281  # N113

$ws.Cells.Item(132, 8).Value = 3318.6099  # H132
$ws.Cells.Item(132, 9).Value = 3179.3215  # I132
$ws.Cells.Item(132, 10).Value = 3618.6155  # J132
$ws.Cells.Item(132, 11).Value = 9537.9645  # K132
$ws.Cells.Item(132, 12).Value = 10855.8465  # L132
$ws.Cells.Item(132, 13).Value = -7007.9645  # M132
$ws.Cells.Item(132, 14).Value = -15915.8465  # N132

$ws.Cells.Item(136, 8).Value = 3095.365  # H136
$ws.Cells.Item(136, 9).Value = 2246.8276  # I136
$ws.Cells.Item(136, 10).Value = 3819.1177  # J136
$ws.Cells.Item(136, 11).Value = 6740.4828  # K136
$ws.Cells.Item(136, 12).Value = 11457.3531  # L136
$ws.Cells.Item(136, 13).Value = -4190.4828  # M136
$ws.Cells.Item(136, 14).Value = -16557.3531  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2090.4883  # H105
$ws.Cells.Item(105, 9).Value = 1931.7646  # I105
$ws.Cells.Item(105, 10).Value = 2690.111  # J105
$ws.Cells.Item(105, 11).Value = 1931.7646  # K105
$ws.Cells.Item(105, 12).Value = 2690.111  # L105
$ws.Cells.Item(105, 13).Value = -184.7646  # M105
$ws.Cells.Item(105, 14).Value = -6184.111  # N105

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 622.7368  # H105
$ws.Cells.Item(105, 9).Value = 602.8  # I105
$ws.Cells.Item(105, 10).Value = 697.5  # J105
$ws.Cells.Item(105, 11).Value = 602.8  # K105
$ws.Cells.Item(105, 12).Value = 697.5  # L105
$ws.Cells.Item(105, 13).Value = 1144.2  # M105
$ws.Cells.Item(105, 14).Value = -4191.5  # N105

$ws.Cells.Item(122, 8).Value = 1198.35  # H122
$ws.Cells.Item(122, 9).Value = 855.5833  # I122
$ws.Cells.Item(122, 10).Value = 1712.5  # J122
$ws.Cells.Item(122, 11).Value = 2566.7499  # K122
$ws.Cells.Item(122, 12).Value = 5137.5  # L122
$ws.Cells.Item(122, 13).Value = -116.7498999999998  # M122
$ws.Cells.Item(122, 14).Value = -10037.5  # N122

$ws.Cells.Item(132, 8).Value = 2245.9  # H132
$ws.Cells.Item(132, 9).Value = 1264.6364  # I132
$ws.Cells.Item(132, 10).Value = 3445.2222  # J132
$ws.Cells.Item(132, 11).Value = 3793.9092  # K132
$ws.Cells.Item(132, 12).Value = 10335.6666  # L132
$ws.Cells.Item(132, 13).Value = -1263.9092  # M132
$ws.Cells.Item(132, 14).Value = -15395.6666  # N132

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 556.5925999999999  # H14
$ws.Cells.Item(14, 9).Value = 556.5925999999999  # I14
$ws.Cells.Item(14, 11).Value = 1669.7778  # K14
$ws.Cells.Item(14, 13).Value = -1496.7778  # M14

$ws.Cells.Item(131, 8).Value = 805.76544  # H131
$ws.Cells.Item(131, 9).Value = 266.69232  # I131
$ws.Cells.Item(131, 10).Value = 908.82355  # J131
$ws.Cells.Item(131, 11).Value = 800.07696  # K131
$ws.Cells.Item(131, 12).Value = 2726.47065  # L131
$ws.Cells.Item(131, 13).Value = 4239.92304  # M131
$ws.Cells.Item(131, 14).Value = -12806.47065  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3883.238  # H102
$ws.Cells.Item(102, 9).Value = 3987.4  # I102
$ws.Cells.Item(102, 11).Value = 3987.4  # K102
$ws.Cells.Item(102, 13).Value = -2365.4  # M102

$ws.Cells.Item(122, 8).Value = 1475.6  # H122
$ws.Cells.Item(122, 9).Value = 1192.6666  # I122
$ws.Cells.Item(122, 11).Value = 3577.9998  # K122
$ws.Cells.Item(122, 13).Value = -1127.9998  # M122

$ws.Cells.Item(132, 8).Value = 4520.909  # H132
$ws.Cells.Item(132, 9).Value = 5324.2354  # I132
$ws.Cells.Item(132, 10).Value = 3667.375  # J132
$ws.Cells.Item(132, 11).Value = 15972.7062  # K132
$ws.Cells.Item(132, 12).Value = 11002.125  # L132
$ws.Cells.Item(132, 13).Value = -13442.7062  # M132
$ws.Cells.Item(132, 14).Value = -16062.125  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2592.3333  # H7
$ws.Cells.Item(7, 9).Value = 2519.9  # I7
$ws.Cells.Item(7, 10).Value = 2658.182  # J7
$ws.Cells.Item(7, 11).Value = 2519.9  # K7
$ws.Cells.Item(7, 12).Value = 2658.182  # L7
$ws.Cells.Item(7, 13).Value = -2407.9  # M7
$ws.Cells.Item(7, 14).Value = -2882.182  # N7

$ws.Cells.Item(16, 8).Value = 1959.56  # H16
$ws.Cells.Item(16, 9).Value = 1934.3  # I16
$ws.Cells.Item(16, 10).Value = 2060.6  # J16
$ws.Cells.Item(16, 11).Value = 1934.3  # K16
$ws.Cells.Item(16, 12).Value = 2060.6  # L16
$ws.Cells.Item(16, 13).Value = -1764.3  # M16
$ws.Cells.Item(16, 14).Value = -2400.6  # N16

$ws.Cells.Item(40, 8).Value = 2383.7334  # H40
$ws.Cells.Item(40, 9).Value = 2318.1538  # I40
$ws.Cells.Item(40, 10).Value = 2810  # J40
$ws.Cells.Item(40, 11).Value = 2318.1538  # K40
$ws.Cells.Item(40, 12).Value = 2810  # L40
$ws.Cells.Item(40, 13).Value = -2182.1538  # M40
$ws.Cells.Item(40, 14).Value = -3082  # N40

$ws.Cells.Item(93, 8).Value = 2119  # H93
$ws.Cells.Item(93, 9).Value = 2300.4443  # I93
$ws.Cells.Item(93, 10).Value = 1885.7142  # J93
$ws.Cells.Item(93, 11).Value = 2300.4443  # K93
$ws.Cells.Item(93, 12).Value = 1885.7142  # L93
$ws.Cells.Item(93, 13).Value = -1052.4443  # M93
$ws.Cells.Item(93, 14).Value = -4381.7142  # N93

$ws.Cells.Item(122, 8).Value = 2769  # H122
$ws.Cells.Item(122, 9).Value = 2600  # I122
$ws.Cells.Item(122, 10).Value = 2966.1667  # J122
$ws.Cells.Item(122, 11).Value = 7800  # K122
$ws.Cells.Item(122, 12).Value = 8898.500100000001  # L122
$ws.Cells.Item(122, 13).Value = -5350  # M122
$ws.Cells.Item(122, 14).Value = -13798.5001  # N122

$ws.Cells.Item(126, 8).Value = 2592.3333  # H126
$ws.Cells.Item(126, 9).Value = 2519.9  # I126
$ws.Cells.Item(126, 10).Value = 2658.182  # J126
$ws.Cells.Item(126, 11).Value = 7559.700000000001  # K126
$ws.Cells.Item(126, 12).Value = 7974.545999999999  # L126
$ws.Cells.Item(126, 13).Value = -5089.700000000001  # M126
$ws.Cells.Item(126, 14).Value = -12914.546  # N126

$ws.Cells.Item(132, 8).Value = 9737  # H132
$ws.Cells.Item(132, 9).Value = 3112.8  # I132
$ws.Cells.Item(132, 11).Value = 9338.400000000001  # K132
$ws.Cells.Item(132, 13).Value = -6808.400000000001  # M132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 17600.5  # H100
$ws.Cells.Item(100, 9).Value = 6800.6665  # I100
$ws.Cells.Item(100, 10).Value = 50000  # J100
$ws.Cells.Item(100, 11).Value = 13601.333  # K100
$ws.Cells.Item(100, 12).Value = 100000  # L100
$ws.Cells.Item(100, 13).Value = -13060.333  # M100
$ws.Cells.Item(100, 14).Value = -101082  # N100

$ws.Cells.Item(108, 8).Value = 0  # H108
$ws.Cells.Item(108, 10).Value = 0  # J108
$ws.Cells.Item(108, 12).ClearContents()  # L108 removed
$ws.Cells.Item(108, 14).Value = 0  # N108

$ws.Cells.Item(113, 8).Value = 315.55554  # H113
$ws.Cells.Item(113, 9).Value = 403.33334  # I113
$ws.Cells.Item(113, 10).Value = 271.66666  # J113
$ws.Cells.Item(113, 11).Value = 1210.00002  # K113
$ws.Cells.Item(113, 12).Value = 814.9999799999999  # L113
$ws.Cells.Item(113, 13).Value = 959.9999800000001  # M113
$ws.Cells.Item(113, 14).Value = -5154.99998  # N113

$ws.Cells.Item(126, 8).Value = 1352.6428  # H126
$ws.Cells.Item(126, 9).Value = 1098.8182  # I126
$ws.Cells.Item(126, 10).Value = 2283.3333  # J126
$ws.Cells.Item(126, 11).Value = 3296.4546  # K126
$ws.Cells.Item(126, 12).Value = 6849.999899999999  # L126
$ws.Cells.Item(126, 13).Value = -826.4546  # M126
$ws.Cells.Item(126, 14).Value = -11789.9999  # N126
